$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.399.56'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '3.416.81'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = "'568.57"
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').Value = "'182.09"
$ws.Range('E6').Value = '  +5.15%  '
$ws.Range('D7').Value = "'0.633"
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('D8').Value = '3.405.88'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('D9').Value = "'1.00"
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +7.09%  '
$ws.Range('D11').Value = "'0.642"
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').Value = "'54.92"
$ws.Range('E12').Value = '  +1.96%  '
$ws.Range('D13').Value = "'0.0000281"
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').Value = "'9.36"
$ws.Range('E14').Value = '  +3.13%  '
$ws.Range('D15').Value = '3.974.33'
$ws.Range('D16').Value = "'18.36"
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.418.84'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = "'0.120"
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = '66.282.47'
$ws.Range('D20').Value = "'12.00"
$ws.Range('E20').Value = '  +1.97%  '
$ws.Range('E21').Value = '  +1.73%  '
$ws.Range('D22').Value = "'466.75"
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = "'14.64"
$ws.Range('E24').Value = '  +9.13%  '
$ws.Range('D25').Value = "'4.15"
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('D26').Value = "'89.80"
$ws.Range('E26').Value = '  +3.73%  '
$ws.Range('D27').Value = "'2.93"
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('D28').Value = "'10.88"
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').Value = "'8.88"
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('D30').Value = "'31.41"
$ws.Range('E30').Value = '  +2.57%  '
$ws.Range('E31').Value = '  +3.95%  '
$ws.Range('D32').Value = "'11.58"
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').Value = "'585.08"
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('D34').Value = "'62.51"
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  +4.79%  '
$ws.Range('D38').Value = "'3.59"
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('D39').Value = "'36.48"
$ws.Range('E39').Value = '  +2.86%  '
$ws.Range('D40').Value = "'0.385"
$ws.Range('E40').Value = '  +4.59%  '
$ws.Range('D41').Value = '0.0₃0759'
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('D42').Value = '3.140.50'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').Value = "'2.95"
$ws.Range('E43').Value = '  +3.84%  '
$ws.Range('D44').Value = "'0.0426"
$ws.Range('E44').Value = '  +2.93%  '
$ws.Range('D45').Value = "'2.54"
$ws.Range('E45').Value = '  +3.49%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = "'2.79"
$ws.Range('E46').Value = '  +18.42%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = "'0.135"
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').Value = "'3.18"
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').Value = "'0.998"
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').Value = "'140.52"
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('D51').Value = "'8.60"
$ws.Range('E51').Value = '  +4.15%  '
